$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

$wsZh.Range("D6").Value = "2016-01-14 04:51:48"
$wsDe.Range("D6").Value = "2016-01-14 04:52:12"
